$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing M:T data (both rows) over to make room for 5 new
# columns (new vCarrierCode..vOffPoint block) - mirrors an Excel
# "Insert Copied/Cut Cells" column-shift so the moved cells keep their
# original types/values (e.g. the text "111.2" in old Q2).
$ws.Range("M1:Q2").Insert(-4161)

# --- Row 3 (new) : core booking columns first ---
$ws.Range("A3").Value = 312
$ws.Range("B3").Value = 1556792
$ws.Range("C3").Value = "SIN"
$ws.Range("D3").Value = "BKK"
$ws.Range("E3").Value = "Test Sprint 9"
$ws.Range("F3").Value = 100
$ws.Range("G3").Value = 1000
$ws.Range("H3").Value = 100
$ws.Range("I3").Value = "AA1"
$ws.Range("J3").Value = "BB2"
$ws.Range("K3").Value = "CC3"
$ws.Range("R3").Value = 112
$ws.Range("S3").Value = 121
$ws.Range("T3").Value = 122
$ws.Range("U3").Value = 100
$ws.Range("V2").Copy($ws.Range("V3"))
$ws.Range("X3").Value = "Sent 1st Booking"

# --- Row 4 (new) : core booking columns first ---
$ws.Range("A4").Value = 312
$ws.Range("B4").Value = 1556792
$ws.Range("C4").Value = "SIN"
$ws.Range("D4").Value = "DPS"
$ws.Range("E4").Value = "Test Sprint 9"
$ws.Range("F4").Value = 100
$ws.Range("G4").Value = 1111
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = "AA1"
$ws.Range("R4").Value = 112
$ws.Range("S4").Value = 121
$ws.Range("T4").Value = 122
$ws.Range("U4").Value = 100
$ws.Range("V2").Copy($ws.Range("V4"))
$ws.Range("X4").Value = "Sent 2nd Booking"

# --- Row 1 header: fill the 5 newly inserted header cells ---
$ws.Range("M1").Value = "vCarrierCode"
$ws.Range("N1").Value = "vFlightNumber"
$ws.Range("O1").Value = "vFlightDate"
$ws.Range("P1").Value = "vBoardPoint"
$ws.Range("Q1").Value = "vOffPoint"

# --- Row 3 / Row 4 new-column data (carrier/flight block) ---
$ws.Range("M3").Value = "ER"
$ws.Range("N3").Value = 111
$ws.Range("O3").Value = 44012
$ws.Range("O3").NumberFormat = "d-mmm-yy"
$ws.Range("P3").Value = "SIN"
$ws.Range("Q3").Value = "BKK"

$ws.Range("M4").Value = "ER"
$ws.Range("N4").Value = 121
$ws.Range("O4").Value = 44012
$ws.Range("O4").NumberFormat = "d-mmm-yy"
$ws.Range("P4").Value = "SIN"
$ws.Range("Q4").Value = "DPS"

# --- Row 1: three brand-new trailing headers U1:W1 ---
$ws.Range("U1").Value = "var_Pieces2"
$ws.Range("V1").Value = "var_Weight2"
$ws.Range("W1").Value = "var_Volume2"

# --- Column widths (bestFit / autofit); columns I:L keep the sheet's
# default width since their short content ("AA1"/"BB2"/"CC3"/"DD4")
# already fits it, so autofit only the columns whose content changed.
$ws.Range("A1:H4").EntireColumn.AutoFit() | Out-Null
$ws.Range("M1:Y4").EntireColumn.AutoFit() | Out-Null

# --- Selection ---
$ws.Range("H18").Select()
